$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.695.40"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "2.490.75"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.06"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.52"
$ws.Range("E6").Value = "  +5.22%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.139"
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.339"
$ws.Range("E11").Value = "  +3.93%  "
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "2.941.01"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.75"
$ws.Range("E14").Value = "  +2.15%  "
$ws.Range("D15").Value = "67.544.75"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").Value = "2.480.16"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.04"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.81"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.66"
$ws.Range("E23").Value = "  +3.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.23"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.23"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D27").Value = "2.616.50"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.994"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "0.0₃0909"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "510.24"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +3.40%  "
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.122"
$ws.Range("E35").Value = "  +7.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.97"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.71"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +3.50%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.330"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.88"
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.44"
$ws.Range("E44").Value = "  +4.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "143.97"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0746"
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.59"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("E51").Value = "  +1.63%  "
